# fix(be): fix and simplify header data for DMS
#
# The title paragraph built "{{ caseId }}" and the trailing line break
# out of extra runs that split the literal text/line-break in odd
# places (leftovers from a manual edit). Merge them back into their
# neighbouring run so the template text reads cleanly as a single run
# per contiguous chunk of formatting, with no functional change to the
# rendered document.

$d = $word.ActiveDocument

# 1) "{{ caseId }" + "}{% if " (two runs) -> "{{ caseId }}{% if " (one run)
$needle1 = "{{ caseId }}{% if "
$found1 = $d.Content.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, $needle1, 2)
if (-not $found1) {
    throw "Could not find the caseId/if text to merge"
}

# 2) " }}){% endif %}" run followed by a separate run that only holds a
#    line break -> merge the break into the "endif" run.
#    A manual line break shows up in Find text as Chr(11) (vertical tab).
$needle2 = " }}){% endif %}" + [char]11
$found2 = $d.Content.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, $needle2, 2)
if (-not $found2) {
    throw "Could not find the endif/line-break text to merge"
}
